$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.703413605690002
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.337289571762085
$ws.Range("D1").Value = 1.514947533607483
$ws.Range("E1").Value = 1.237128615379333
